$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are text (inlineStr) even though many look numeric or
# percentage-like (e.g. "309.36", "-3.58%"). Force each target range to Text
# format before writing so Excel does not auto-convert them into numbers /
# percentages.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '309.36'
$ws.Range("E2").Value = '-3.58%'

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '50.33'
$ws.Range("E3").Value = '2.95%'

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '5.194'
$ws.Range("E4").Value = '-1.15%'

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07773'
$ws.Range("E5").Value = '-4.24%'

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '4.495'
$ws.Range("E6").Value = '-2.07%'

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '1.341'
$ws.Range("E7").Value = '10.92%'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-5.01%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1211'

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1982'
$ws.Range("E10").Value = '1.78%'

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04799'
$ws.Range("E11").Value = '4.83%'

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09437'
$ws.Range("E12").Value = '0.11%'

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1044'
$ws.Range("E13").Value = '-0.72%'

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001269'
$ws.Range("E14").Value = '-4.90%'

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005792'
$ws.Range("E15").Value = '-1.37%'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,015.68%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.335'

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '2.433'
$ws.Range("E18").Value = '0.33%'

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3476'
$ws.Range("E19").Value = '1.70%'

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '8.003'
$ws.Range("E20").Value = '-1.10%'

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1362'
$ws.Range("E21").Value = '-0.57%'

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '0.3092'
$ws.Range("E22").Value = '-1.07%'

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04162'
$ws.Range("E23").Value = '-0.11%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.79%'

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003947'
$ws.Range("E25").Value = '-7.26%'

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001350'
$ws.Range("E26").Value = '-0.05%'

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02606'
$ws.Range("E38").Value = '-3.96%'

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06108'
$ws.Range("E39").Value = '6.37%'

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01100'
$ws.Range("E40").Value = '74.49%'

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007950'
$ws.Range("E41").Value = '2.61%'

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1422'
$ws.Range("E42").Value = '-1.46%'

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008372'
$ws.Range("E43").Value = '8.86%'

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008356'
$ws.Range("E44").Value = '3.15%'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '5.43%'

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007226'
$ws.Range("E46").Value = '3.28%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.06%'

$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '0.002618'
$ws.Range("E48").Value = '-34.59%'

$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '0.05318'
$ws.Range("E49").Value = '-15.49%'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.06%'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.06%'
